$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 (styled like the other header cells, s="1")
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data for columns I (col 9) and J (col 10), rows 2..68
$values = @(
    @(5, 6),
    @(6, 7),
    @(4, 5),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(5, 6),
    @(6, 6),
    @(5, 5),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(6, 7),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(7, 9),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(6, 6),
    @(6, 7),
    @(4, 5),
    @(6, 6),
    @(7, 8),
    @(9, 9),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(5, 6),
    @(5, 6),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(6, 7),
    @(7, 7),
    @(6, 7),
    @(7, 8),
    @(6, 7),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(1, 1),
    @(8, 8),
    @(7, 7),
    @(3, 3)
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
